$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New forecast-error values (Component Analysis filtering applied)
$values = @(
    @(0.08102716014508599, 0.1849562581409688,  0.05080724786734196, 0.225404631423895,  0.2177199557413915, 15),
    @(0.3284388860093228,  0.3436510707338022,  0.1900396029840986,  0.4359353197254136, 0.2974694575110258, 14),
    @(0.5008085514463395,  0.5008085514463395,  0.4390583335671841,  0.6626147701094386, 0.4515933521755904, 13),
    @(0.6271289531314115,  0.6271289531314115,  0.6001355647093024,  0.77468417094278,   0.4750251754326699, 12),
    @(0.6196096897755812,  0.6196096897755812,  0.4847206317237406,  0.6962188102340676, 0.3329938595018584, 11),
    @(0.4878350037768636,  0.4878350037768636,  0.2918424992789185,  0.5402244897067501, 0.2446301252661204, 10),
    @(0.408887426949955,   0.408887426949955,   0.2151763681766355,  0.4638710684841592, 0.2323485964907919, 9),
    @(0.399669795446622,   0.399669795446622,   0.2150887163119547,  0.4637765801676004, 0.2577272300389153, 6),
    @(0.2327354214815493,  0.2327354214815493,  0.07574075157082642, 0.2752103769315874, 0.1798956996093793, 3),
    @(0.2378331735378737,  0.2378331735378737,  0.05656461843509633, 0.2378331735378737, $null,              1)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $vals = $values[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    if ($null -eq $vals[4]) {
        $ws.Cells.Item($row, 6).ClearContents()
    } else {
        $ws.Cells.Item($row, 6).Value = $vals[4]
    }
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
